$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of G1 (header) onto H1, then set the value
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
